$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Plan las pruebas" -> "Plan de" + " pruebas" (two runs)
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Plan las pruebas", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$start1 = $rng1.Start

# Replace "las" with "de" in place (keeps it within the same run for now).
$rLas = $d.Range($start1 + 5, $start1 + 8)
$rLas.Text = "de"

# Force a run boundary right after "Plan de" by nudging the font of the
# first chunk; this keeps " pruebas" as an untouched run carrying the
# original, full run formatting (rFonts/color/lang all preserved), while
# "Plan de" becomes its own run.
$rFirst = $d.Range($start1, $start1 + 7)
$rFirst.Font.Name = "Calibri"

# ---------------------------------------------------------------------
# Change 2: "Documento de Evaluación de Pruebas y Establecimiento de
#            Métricas" -> "Documento de Evaluación de Pruebas " (keeps
#            the trailing space, single run).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Documento de Evaluación de Pruebas y Establecimiento de Métricas", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Documento de Evaluación de Pruebas ", 2)
